# Generate Report for Handback
# The 8e42f056-ec0a-4c7f-a146-0e3fc7aff847 file has now been handed back
# (target generated + handed back to the source repo), so it moves to the
# top of the status report, its status flips to "Handed back: in sync with
# en-US", and its per-language rows grow a populated Latest Target
# File / Latest Handback File / Latest Handback DateTime trio. The
# 38f23d60-677f-49ef-a9ec-0d49dee3e928 file is unaffected content-wise but
# drops to the second row.

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob"
$xlfBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob"

$hash8e   = "3e3c212bf10f6a61559469591afa0cbdd0613bd9"
$hash38f  = "62a685063ad850ba21bc8b689f56889bb6379ba3"

$md8e   = "$mdBase/$hash8e/e2e/8e42f056-ec0a-4c7f-a146-0e3fc7aff847.md"
$md38f  = "$mdBase/$hash38f/e2e/38f23d60-677f-49ef-a9ec-0d49dee3e928.md"

$xlf8eZh  = "$xlfBase/3def75a57439140d18d4a8a1fd01479e10ea07b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e42f056-ec0a-4c7f-a146-0e3fc7aff847.86da594c26610804b305b7406ed3306746770926.zh-cn.xlf"
$xlf38fZh = "$xlfBase/6fd77e66b5964f9a7f899c4099aceb11b0e5c6f3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/38f23d60-677f-49ef-a9ec-0d49dee3e928.e0be370d9f7de06c002332827821f0a40ee39c34.zh-cn.xlf"

$xlf8eDe  = "$xlfBase/1ca63907d52d775064435229e82968bc54a6d8f5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e42f056-ec0a-4c7f-a146-0e3fc7aff847.86da594c26610804b305b7406ed3306746770926.de-de.xlf"
$xlf38fDe = "$xlfBase/82d856fdb7423131c55eef42db4e40afd1a845f2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/38f23d60-677f-49ef-a9ec-0d49dee3e928.e0be370d9f7de06c002332827821f0a40ee39c34.de-de.xlf"

$name8e  = "8e42f056-ec0a-4c7f-a146-0e3fc7aff847.md"
$name38f = "38f23d60-677f-49ef-a9ec-0d49dee3e928.md"
$xlf8eZhName  = "8e42f056-ec0a-4c7f-a146-0e3fc7aff847.86da594c26610804b305b7406ed3306746770926.zh-cn.xlf"
$xlf38fZhName = "38f23d60-677f-49ef-a9ec-0d49dee3e928.e0be370d9f7de06c002332827821f0a40ee39c34.zh-cn.xlf"
$xlf8eDeName  = "8e42f056-ec0a-4c7f-a146-0e3fc7aff847.86da594c26610804b305b7406ed3306746770926.de-de.xlf"
$xlf38fDeName = "38f23d60-677f-49ef-a9ec-0d49dee3e928.e0be370d9f7de06c002332827821f0a40ee39c34.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview" - swap the two rows so the freshly handed-back file
# (8e42f056) sits on row 2, and update its status / date.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = $name8e
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-27-18 10:27:12"

$ov.Range("A3").Value = $name38f
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-26-18 10:26:51"

$ov.Hyperlinks.Add($ov.Range("A2"), $md8e, "", "", $name8e)
$ov.Hyperlinks.Add($ov.Range("A3"), $md38f, "", "", $name38f)

# ---------------------------------------------------------------------
# Sheet "zh-cn" - same row reorder; the handed-back row also gains its
# Latest Target File / Latest Handback File / Latest Handback DateTime.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $name8e
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = $xlf8eZhName
$zh.Range("E2").Value = "2016-03-18 10:27:10"
$zh.Range("F2").Value = $name8e
$zh.Range("G2").Value = $xlf8eZhName
$zh.Range("H2").Value = "2016-03-18 10:27:27"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = $name38f
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = $xlf38fZhName
$zh.Range("E3").Value = "2016-03-18 10:26:48"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A2"), $md8e, "", "", $name8e)
$zh.Hyperlinks.Add($zh.Range("B2"), $md8e, "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), $xlf8eZh, "", "", $xlf8eZhName)
$zh.Hyperlinks.Add($zh.Range("F2"), $md8e, "", "", $name8e)
$zh.Hyperlinks.Add($zh.Range("G2"), $xlf8eZh, "", "", $xlf8eZhName)

$zh.Hyperlinks.Add($zh.Range("A3"), $md38f, "", "", $name38f)
$zh.Hyperlinks.Add($zh.Range("B3"), $md38f, "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D3"), $xlf38fZh, "", "", $xlf38fZhName)

# ---------------------------------------------------------------------
# Sheet "de-de" - mirror of the zh-cn changes, de-de file names/dates.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A2").Value = $name8e
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = $xlf8eDeName
$de.Range("E2").Value = "2016-03-18 10:27:12"
$de.Range("F2").Value = $name8e
$de.Range("G2").Value = $xlf8eDeName
$de.Range("H2").Value = "2016-03-18 10:27:32"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = $name38f
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = $xlf38fDeName
$de.Range("E3").Value = "2016-03-18 10:26:51"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Add($de.Range("A2"), $md8e, "", "", $name8e)
$de.Hyperlinks.Add($de.Range("B2"), $md8e, "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), $xlf8eDe, "", "", $xlf8eDeName)
$de.Hyperlinks.Add($de.Range("F2"), $md8e, "", "", $name8e)
$de.Hyperlinks.Add($de.Range("G2"), $xlf8eDe, "", "", $xlf8eDeName)

$de.Hyperlinks.Add($de.Range("A3"), $md38f, "", "", $name38f)
$de.Hyperlinks.Add($de.Range("B3"), $md38f, "", "", ".md")
$de.Hyperlinks.Add($de.Range("D3"), $xlf38fDe, "", "", $xlf38fDeName)

Write-Output "Handback report regenerated"
